# Creator_Scout.xlsx update: "added the new data"
# - Adds a new "apify filter" (Yes/No) column D to the Instructions sheet
# - Expands several of the explanatory "Note" texts in column C
# - Makes the Instructions sheet the active / selected tab (was Sheet1)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# ---- New column D header -------------------------------------------------
$ws.Range("D1").Value = "apify filter"
$ws.Range("D1").Interior.Color = 65535   # reuse the same yellow header fill as A1:C1

# ---- Row 2 : Description_Keywords ----------------------------------------
$ws.Range("B2").Value = "String"
$ws.Range("C2").Value = "String`twords must be comma seperated else will be considered as one,  Your keyword will apply user profiles."
$ws.Range("D2").Value = "Yes"
$ws.Rows.Item(2).RowHeight = 43.5

# ---- Row 3 : Category ------------------------------------------------------
$ws.Range("B3").Value = "String"
$ws.Range("C3").Value = "String words must be comma seperated in case of multiple else will be considered as one ,  Your keyword will apply to both videos description and user profiles."
$ws.Range("D3").Value = "Yes"
$ws.Rows.Item(3).RowHeight = 57.75

# ---- Row 4 : Language -------------------------------------------------------
$ws.Range("B4").Value = "String"
$ws.Range("C4").Value = "proper language words must be insterted ,Only 1 language be provided at a time. eg-> German/ English"
$ws.Range("D4").Value = "No"
$ws.Rows.Item(4).RowHeight = 43.5

# ---- Row 5 : Last_Post_Period(7,14,30) --------------------------------------
$ws.Range("B5").Value = "number"
$ws.Range("C5").Value = "days cout in number , if 0 then will not be considered ,eg-> 5 , Only videos uploaded after the days of how old the scraped videos should be in days. Putting 1 will get you only today's posts, 2 - yesterday's and today's, and so on."
$ws.Range("D5").Value = "Yes"
$ws.Rows.Item(5).RowHeight = 72.75

# ---- Row 6 : Min_Followers ---------------------------------------------------
$ws.Range("B6").Value = "number"
$ws.Range("C6").Value = "Min. follower in number , if 0 then will not be considered ,eg-> 5 , Scrapes only profiles with more then  number on followers in the profile"
$ws.Range("D6").Value = "Yes"
$ws.Rows.Item(6).RowHeight = 43.5

# ---- Row 7 : Max_Followers ----------------------------------------------------
$ws.Range("B7").Value = "number"
$ws.Range("C7").Value = "Max. follower in number , if 0 then will not be considered ,eg-> 5,  Scrapes only profiles with more then  number on followers in the profile"
$ws.Range("D7").Value = "Yes"
$ws.Rows.Item(7).RowHeight = 43.5

# ---- Row 8 : Min_Average_Likes -------------------------------------------------
$ws.Range("B8").Value = "number"
$ws.Range("C8").Value = "Min. Avg. Like in number , if 0 then will not be considered ,eg-> 5 , will calculate at our side that average like for. the profile and filter out the use data based on the minimum average data . "
$ws.Range("D8").Value = "No"
$ws.Rows.Item(8).RowHeight = 57.75

# ---- Row 9 : Min_Avg_Comment -----------------------------------------------------
$ws.Range("B9").Value = "number"
$ws.Range("C9").Value = "Min. Avg Comment follower in number , if 0 then will not be considered ,eg-> 5 ,will calculate at our side that average comment  for the profile and filter out the use data based on the minimum average data "
$ws.Range("D9").Value = "No"
$ws.Rows.Item(9).RowHeight = 57.75

# ---- Row 10 : Min_Engagement_Rate -------------------------------------------------
$ws.Range("B10").Value = "number"
$ws.Range("C10").Value = "post having engagement rate % greater then provided number will be considered for result, The value provided will filtered at server and profile having the min. engagement greater then number will be provided."
$ws.Range("D10").Value = "No"
$ws.Rows.Item(10).RowHeight = 72.75

# ---- Row 11 : Number_of_Required_Results --------------------------------------------
$ws.Range("B11").Value = "number"
$ws.Range("C11").Value = "the number of requred data to fetch from the apify after filter (internal filters are not included in this.) "
$ws.Range("D11").Value = "Yes"
$ws.Rows.Item(11).RowHeight = 29.25

# ---- Row 12 : Location --------------------------------------------------------------
$ws.Range("B12").Value = "String"
$ws.Range("C12").Value = "country string if need result from specific country , else will be on empty no country filter will be applied. "
$ws.Range("D12").Value = "Yes"
$ws.Rows.Item(12).RowHeight = 29.25

# ---- Row 13 : Contact_Info_Available --------------------------------------------------
$ws.Range("B13").Value = "String ('yes' / 'No' / Empty)"
$ws.Range("C13").Value = "if yes then only return the result of the user that as contactable info. in case of  no and empty this filter will not be used. "
$ws.Range("D13").Value = "No"
$ws.Rows.Item(13).RowHeight = 43.5

# ---- Make "Instructions" the active / selected sheet (was "Sheet1") -------
$ws.Activate() | Out-Null
$ws.Range("D1").Select() | Out-Null
